$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 120
$ws1.Range("F7").Value = 11994
$ws1.Range("F8").Value = 4456
$ws1.Range("F11").Value = 29
$ws1.Range("F17").Value = 5195
$ws1.Range("F21").Value = 11402
$ws1.Range("F22").Value = 11424

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 120
$ws4.Range("F7").Value = 11994
$ws4.Range("F8").Value = 4456
$ws4.Range("F11").Value = 29
$ws4.Range("F18").Value = 5195
$ws4.Range("F22").Value = 11402
$ws4.Range("F23").Value = 11424
